$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" (column C) values to reflect the new run's log data.
# The new run converges faster, so the recorded best-fitness-so-far values
# change for generations 0-61 (rows 2-63); later generations already matched
# the new value (7293) and remain unchanged.
$ws.Range("C2:C28").Value = 7928
$ws.Range("C29:C30").Value = 7916
$ws.Range("C31:C32").Value = 7632
$ws.Range("C33:C63").Value = 7293
